# Update the timestamp embedded in the test e-mail addresses on the
# "UsuariosRegistro" sheet (column C, rows 2-6) from 20251109_004215 to
# 20251109_005042. Since "LoginData" sheet A2/A3 share the same strings
# (juan.perez / maria.gonzalez addresses), updating the text in place
# keeps every usage of those addresses consistent.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_004215"
$newStamp = "20251109_005042"

$ws1 = $wb.Worksheets.Item("UsuariosRegistro")
$regCells = @("C2", "C3", "C4", "C5", "C6")
foreach ($addr in $regCells) {
    $cell = $ws1.Range($addr)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}

# "LoginData" reuses the juan.perez and maria.gonzalez addresses in A2:A3,
# so update those cells as well to keep every occurrence in sync.
$ws2 = $wb.Worksheets.Item("LoginData")
$loginCells = @("A2", "A3")
foreach ($addr in $loginCells) {
    $cell = $ws2.Range($addr)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}
